# Updated cryptos list on Mon Oct  7 17:53:43 UTC 2024 with GitHub Actions
#
# This script refreshes the "Price" (column D) and "Volume(1h)" (column E)
# figures scraped from coinranking.com, and re-sorts two pairs of rows whose
# relative ranking flipped since the last run (InternetComputer(DFINITY) /
# Binance-PegBSC-USD at rows 31-32, and Stacks / Monero at rows 42-43).
#
# Price/volume cells are stored as plain TEXT in the sheet (e.g. "63.727.29",
# "1.00", "  +1.74%  ") rather than numbers, so we force text entry for any
# replacement value that would otherwise be auto-parsed into a numeric type
# by Excel (losing formatting like trailing zeros, or multi-dot separators).

function Set-CellText {
    param($ws, [string]$addr, [string]$val)
    $cell = $ws.Range($addr)
    if ($val -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Would be auto-coerced to a number by Excel - force text with a
        # leading quote-prefix, then strip the resulting cell style back to
        # Normal so no visible/format change is introduced.
        $cell.Value = "'" + $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "63.727.29"
Set-CellText $ws "E2" "  +1.74%  "
Set-CellText $ws "D3" "2.476.47"
Set-CellText $ws "D4" "1.00"
Set-CellText $ws "E4" "  -0.09%  "
Set-CellText $ws "D5" "576.36"
Set-CellText $ws "E5" "  +1.74%  "
Set-CellText $ws "D6" "149.08"
Set-CellText $ws "E6" "  +2.70%  "
Set-CellText $ws "E7" "  +0.02%  "
Set-CellText $ws "E8" "  +2.00%  "
Set-CellText $ws "D9" "2.474.44"
Set-CellText $ws "E9" "  +1.45%  "
Set-CellText $ws "E10" "  +1.22%  "
Set-CellText $ws "D11" "0.157"
Set-CellText $ws "E11" "  +1.27%  "
Set-CellText $ws "D12" "5.31"
Set-CellText $ws "E12" "  +1.15%  "
Set-CellText $ws "E13" "  +1.60%  "
Set-CellText $ws "D14" "27.26"
Set-CellText $ws "E14" "  +1.63%  "
Set-CellText $ws "E15" "  -0.01%  "
Set-CellText $ws "D16" "2.923.43"
Set-CellText $ws "E16" "  +1.54%  "
Set-CellText $ws "D17" "63.492.96"
Set-CellText $ws "E17" "  +1.45%  "
Set-CellText $ws "D18" "2.479.31"
Set-CellText $ws "E18" "  +1.87%  "
Set-CellText $ws "E19" "  +2.36%  "
Set-CellText $ws "D20" "7.42"
Set-CellText $ws "E20" "  +7.32%  "
Set-CellText $ws "D21" "331.82"
Set-CellText $ws "E21" "  +2.42%  "
Set-CellText $ws "D22" "4.24"
Set-CellText $ws "E22" "  +1.90%  "
Set-CellText $ws "D23" "2.12"
Set-CellText $ws "E23" "  +19.13%  "
Set-CellText $ws "E24" "  +0.36%  "
Set-CellText $ws "D25" "66.03"
Set-CellText $ws "E25" "  -1.73%  "
Set-CellText $ws "D26" "631.01"
Set-CellText $ws "E26" "  +11.46%  "
Set-CellText $ws "D27" "9.25"
Set-CellText $ws "E27" "  +6.60%  "
Set-CellText $ws "D28" "0.0000106"
Set-CellText $ws "E28" "  +4.69%  "
Set-CellText $ws "D29" "1.55"
Set-CellText $ws "E29" "  +7.08%  "
Set-CellText $ws "D30" "2.603.38"
Set-CellText $ws "B31" "InternetComputer(DFINITY)"
Set-CellText $ws "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText $ws "D31" "8.44"
Set-CellText $ws "E31" "  +0.75%  "
Set-CellText $ws "B32" "Binance-PegBSC-USD"
Set-CellText $ws "C32" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-CellText $ws "D32" "0.998"
Set-CellText $ws "E32" "  -0.04%  "
Set-CellText $ws "E33" "  -2.08%  "
Set-CellText $ws "E34" "  +2.51%  "
Set-CellText $ws "D35" "5.27"
Set-CellText $ws "E35" "  +8.64%  "
Set-CellText $ws "E36" "  +1.18%  "
Set-CellText $ws "D37" "1.00"
Set-CellText $ws "E37" "  +0.11%  "
Set-CellText $ws "E38" "  +0.30%  "
Set-CellText $ws "D39" "5.54"
Set-CellText $ws "E39" "  +2.06%  "
Set-CellText $ws "D40" "18.94"
Set-CellText $ws "E40" "  +0.82%  "
Set-CellText $ws "E41" "  +14.65%  "
Set-CellText $ws "B42" "Stacks"
Set-CellText $ws "C42" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-CellText $ws "D42" "1.82"
Set-CellText $ws "E42" "  +0.31%  "
Set-CellText $ws "B43" "Monero"
Set-CellText $ws "C43" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText $ws "D43" "147.27"
Set-CellText $ws "E43" "  -0.69%  "
Set-CellText $ws "E44" "  -0.38%  "
Set-CellText $ws "D45" "151.14"
Set-CellText $ws "E45" "  +2.09%  "
Set-CellText $ws "D46" "3.80"
Set-CellText $ws "E46" "  +3.52%  "
Set-CellText $ws "D47" "21.61"
Set-CellText $ws "E47" "  +5.73%  "
Set-CellText $ws "E48" "  +1.19%  "
Set-CellText $ws "E49" "  +1.35%  "
Set-CellText $ws "E50" "  +2.99%  "
Set-CellText $ws "D51" "0.0922"
Set-CellText $ws "E51" "  -0.31%  "
